$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 World Cup")

# Results for the last group-stage matches (Jun 28, 2018) — Groups G & H
# Row 51: Inglaterra vs Bélgica -> 0 - 1
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1

# Row 52: Panamá vs Túnez -> 1 - 2
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 2

# Row 53: Japón vs Polonia -> 0 - 1
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 1

# Row 54: Senegal vs Colombia -> 0 - 1
$ws.Range("F54").Value = 0
$ws.Range("G54").Value = 1

# Update the view so the active selection/scroll position reflects the edit
$ws.Activate()
$excel.Goto($ws.Range("Q10"), $true)
$ws.Range("F52").Select()
